$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 value from "Australia" to "Germany" (new shared string entry)
$ws.Range("D2").Value = "Germany"

# Move/select the active cell to E18 as recorded in the saved view state
$ws.Range("E18").Select()
